$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 461.35715
$ws.Range("I2").Value = 391.72726
$ws.Range("J2").Value = 716.6667
$ws.Range("K2").Value = 391.72726
$ws.Range("L2").Value = 716.6667
$ws.Range("M2").Value = -278.72726
$ws.Range("N2").Value = -942.6667
$ws.Range("H17").Value = 1761.0714
$ws.Range("I17").Value = 1444
$ws.Range("J17").Value = 1772.8148
$ws.Range("K17").Value = 4332
$ws.Range("L17").Value = 5318.4444
$ws.Range("M17").Value = -4164
$ws.Range("N17").Value = -5654.4444
$ws.Range("H38").Value = 1315.1538
$ws.Range("I38").Value = 144.44444
$ws.Range("J38").Value = 3949.25
$ws.Range("K38").Value = 433.33332
$ws.Range("L38").Value = 11847.75
$ws.Range("M38").Value = -61.33331999999996
$ws.Range("N38").Value = -12591.75
$ws.Range("H40").Value = 1722.1538
$ws.Range("I40").Value = 1712.375
$ws.Range("J40").Value = 1737.8
$ws.Range("K40").Value = 1712.375
$ws.Range("L40").Value = 1737.8
$ws.Range("M40").Value = -1537.375
$ws.Range("N40").Value = -2087.8
$ws.Range("H100").Value = 2021.6316
$ws.Range("I100").Value = 1801
$ws.Range("K100").Value = 1801
$ws.Range("M100").Value = -1260
$ws.Range("H106").Value = 9847.1
$ws.Range("I106").Value = 3789.7693
$ws.Range("K106").Value = 3789.7693
$ws.Range("M106").Value = -3158.7693
$ws.Range("H115").Value = 1505.625
$ws.Range("I115").Value = 1623.6666
$ws.Range("K115").Value = 4870.9998
$ws.Range("M115").Value = -3303.9998
$ws.Range("H116").Value = 30027.555
$ws.Range("I116").Value = 23392.572
$ws.Range("K116").Value = 23392.572
$ws.Range("M116").Value = -19950.572
$ws.Range("H118").Value = 379
$ws.Range("I118").Value = 369
$ws.Range("J118").Value = 419
$ws.Range("K118").Value = 1107
$ws.Range("L118").Value = 1257
$ws.Range("M118").Value = 550
$ws.Range("N118").Value = -4571
$ws.Range("H132").Value = 10817.2
$ws.Range("I132").Value = 10999.625
$ws.Range("J132").Value = 6439
$ws.Range("K132").Value = 32998.875
$ws.Range("L132").Value = 19317
$ws.Range("M132").Value = -30468.875
$ws.Range("N132").Value = -24377
$ws.Range("H137").Value = 4419.222
$ws.Range("I137").Value = 2310.75
$ws.Range("J137").Value = 8086.1304
$ws.Range("K137").Value = 6932.25
$ws.Range("L137").Value = 24258.3912
$ws.Range("M137").Value = -4382.25
$ws.Range("N137").Value = -29358.3912
$ws.Range("H138").Value = 5039.393
$ws.Range("I138").Value = 5164.12
$ws.Range("J138").Value = 4000
$ws.Range("K138").Value = 15492.36
$ws.Range("L138").Value = 12000
$ws.Range("M138").Value = -10352.36
$ws.Range("N138").Value = -22280

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 130002.69
$ws.Range("I32").Value = 151515.7
$ws.Range("K32").Value = 151515.7
$ws.Range("M32").Value = -151228.7
$ws.Range("H45").Value = 2504
$ws.Range("I45").Value = 1727.1111
$ws.Range("K45").Value = 1727.1111
$ws.Range("M45").Value = -1350.1111
$ws.Range("H61").Value = 1902.6342
$ws.Range("I61").Value = 1887.7693
$ws.Range("J61").Value = 2192.5
$ws.Range("K61").Value = 1887.7693
$ws.Range("L61").Value = 2192.5
$ws.Range("M61").Value = -1675.7693
$ws.Range("N61").Value = -2616.5
$ws.Range("H74").Value = 4853.6763
$ws.Range("I74").Value = 1395.3462
$ws.Range("J74").Value = 16093.25
$ws.Range("K74").Value = 1395.3462
$ws.Range("L74").Value = 16093.25
$ws.Range("M74").Value = -521.3462
$ws.Range("N74").Value = -17841.25
$ws.Range("H77").Value = 4853.6763
$ws.Range("I77").Value = 1395.3462
$ws.Range("J77").Value = 16093.25
$ws.Range("K77").Value = 6976.731
$ws.Range("L77").Value = 80466.25
$ws.Range("M77").Value = -2608.731
$ws.Range("N77").Value = -89202.25
$ws.Range("H88").Value = 1410.8518
$ws.Range("I88").Value = 619.5
$ws.Range("J88").Value = 1744.0526
$ws.Range("K88").Value = 619.5
$ws.Range("L88").Value = 1744.0526
$ws.Range("M88").Value = -213.5
$ws.Range("N88").Value = -2556.0526
$ws.Range("H91").Value = 1410.8518
$ws.Range("I91").Value = 619.5
$ws.Range("J91").Value = 1744.0526
$ws.Range("K91").Value = 619.5
$ws.Range("L91").Value = 1744.0526
$ws.Range("M91").Value = 784.5
$ws.Range("N91").Value = -4552.0526
$ws.Range("H102").Value = 6456.6665
$ws.Range("I102").Value = 5216.6665
$ws.Range("K102").Value = 5216.6665
$ws.Range("M102").Value = -3594.6665
$ws.Range("H110").Value = 769.8
$ws.Range("I110").Value = 769.8
$ws.Range("K110").Value = 769.8
$ws.Range("M110").Value = 1275.2
$ws.Range("H132").Value = 2086821.1
$ws.Range("I132").Value = 2780995
$ws.Range("J132").Value = 4299.3335
$ws.Range("K132").Value = 8342985
$ws.Range("L132").Value = 12898.0005
$ws.Range("M132").Value = -8340455
$ws.Range("N132").Value = -17958.0005
$ws.Range("H136").Value = 1902.6342
$ws.Range("I136").Value = 1887.7693
$ws.Range("J136").Value = 2192.5
$ws.Range("K136").Value = 5663.3079
$ws.Range("L136").Value = 6577.5
$ws.Range("M136").Value = -3113.3079
$ws.Range("N136").Value = -11677.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 43537
$ws.Range("I20").Value = 60068.832
$ws.Range("J20").Value = 1026.5714
$ws.Range("K20").Value = 60068.832
$ws.Range("L20").Value = 1026.5714
$ws.Range("M20").Value = -59821.832
$ws.Range("N20").Value = -1520.5714
$ws.Range("H44").Value = 9825
$ws.Range("H86").Value = 2630.1333
$ws.Range("I86").Value = 2700.3635
$ws.Range("K86").Value = 2700.3635
$ws.Range("M86").Value = -1577.3635
$ws.Range("H89").Value = 2630.1333
$ws.Range("I89").Value = 2700.3635
$ws.Range("K89").Value = 13501.8175
$ws.Range("M89").Value = -7885.817499999999
$ws.Range("H99").Value = 7351.6113
$ws.Range("I99").Value = 7725.2354
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 7725.2354
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = -6227.2354
$ws.Range("N99").Value = -3996
$ws.Range("H105").Value = 3845.5
$ws.Range("I105").Value = 1507.6
$ws.Range("J105").Value = 5144.3335
$ws.Range("K105").Value = 1507.6
$ws.Range("L105").Value = 5144.3335
$ws.Range("M105").Value = 239.4000000000001
$ws.Range("N105").Value = -8638.333500000001
$ws.Range("H134").Value = 4691.923
$ws.Range("I134").Value = 2107.7026
$ws.Range("K134").Value = 6323.1078
$ws.Range("M134").Value = -3788.1078

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 45584.434
$ws.Range("I16").Value = 2153.2856
$ws.Range("J16").Value = 113144
$ws.Range("K16").Value = 2153.2856
$ws.Range("L16").Value = 113144
$ws.Range("M16").Value = -1866.2856
$ws.Range("N16").Value = -113718
$ws.Range("H31").Value = 3536.1191
$ws.Range("I31").Value = 3616.9
$ws.Range("J31").Value = 3334.1667
$ws.Range("K31").Value = 3616.9
$ws.Range("L31").Value = 3334.1667
$ws.Range("M31").Value = -3321.9
$ws.Range("N31").Value = -3924.1667
$ws.Range("H34").Value = 3536.1191
$ws.Range("I34").Value = 3616.9
$ws.Range("J34").Value = 3334.1667
$ws.Range("K34").Value = 3616.9
$ws.Range("L34").Value = 3334.1667
$ws.Range("M34").Value = -3414.9
$ws.Range("N34").Value = -3738.1667
$ws.Range("H36").Value = 14998
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H40").Value = 14998
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H54").Value = 36496
$ws.Range("J54").Value = 36496
$ws.Range("L54").Value = 36496
$ws.Range("N54").Value = -37812
$ws.Range("H58").Value = 6301.1724
$ws.Range("I58").Value = 2321.1428
$ws.Range("J58").Value = 16748.75
$ws.Range("K58").Value = 2321.1428
$ws.Range("L58").Value = 16748.75
$ws.Range("M58").Value = -2118.1428
$ws.Range("N58").Value = -17154.75
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H96").Value = 12750
$ws.Range("J96").Value = 12750
$ws.Range("L96").Value = 12750
$ws.Range("N96").Value = -18242
$ws.Range("H113").Value = 45584.434
$ws.Range("I113").Value = 2153.2856
$ws.Range("J113").Value = 113144
$ws.Range("K113").Value = 2153.2856
$ws.Range("L113").Value = 113144
$ws.Range("M113").Value = 16.71439999999984
$ws.Range("N113").Value = -117484
$ws.Range("H132").Value = 2050.1843
$ws.Range("I132").Value = 1842.7812
$ws.Range("J132").Value = 3156.3333
$ws.Range("K132").Value = 5528.3436
$ws.Range("L132").Value = 9468.999899999999
$ws.Range("M132").Value = -2998.3436
$ws.Range("N132").Value = -14528.9999
$ws.Range("H134").Value = 2060.889
$ws.Range("I134").Value = 1693.5
$ws.Range("K134").Value = 5080.5
$ws.Range("M134").Value = -2545.5
$ws.Range("H136").Value = 6301.1724
$ws.Range("I136").Value = 2321.1428
$ws.Range("J136").Value = 16748.75
$ws.Range("K136").Value = 6963.428400000001
$ws.Range("L136").Value = 50246.25
$ws.Range("M136").Value = -4413.428400000001
$ws.Range("N136").Value = -55346.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 453.72223
$ws.Range("I23").Value = 238.8
$ws.Range("J23").Value = 536.38464
$ws.Range("K23").Value = 716.4000000000001
$ws.Range("L23").Value = 1609.15392
$ws.Range("M23").Value = -481.4000000000001
$ws.Range("N23").Value = -2079.15392
$ws.Range("H32").Value = 3798.4443
$ws.Range("J32").Value = 3966
$ws.Range("L32").Value = 11898
$ws.Range("N32").Value = -12464
$ws.Range("H70").Value = 1583.1666
$ws.Range("I70").Value = 1299.8
$ws.Range("K70").Value = 3899.4
$ws.Range("M70").Value = -3584.4
$ws.Range("H73").Value = 1583.1666
$ws.Range("I73").Value = 1299.8
$ws.Range("K73").Value = 3899.4
$ws.Range("M73").Value = -2807.4
$ws.Range("H75").Value = 1664.6666
$ws.Range("I75").Value = 496.33334
$ws.Range("K75").Value = 1489.00002
$ws.Range("M75").Value = -491.0000199999999
$ws.Range("H78").Value = 1664.6666
$ws.Range("I78").Value = 496.33334
$ws.Range("K78").Value = 4467.00006
$ws.Range("M78").Value = 524.9999399999997
$ws.Range("H107").Value = 536
$ws.Range("J107").Value = 731.5
$ws.Range("L107").Value = 2194.5
$ws.Range("N107").Value = -6034.5
$ws.Range("H121").Value = 22729372
$ws.Range("I121").Value = 439.6
$ws.Range("J121").Value = 29414352
$ws.Range("K121").Value = 1318.8
$ws.Range("L121").Value = 88243056
$ws.Range("M121").Value = -8.800000000000182
$ws.Range("N121").Value = -88245676
$ws.Range("H129").Value = 4661233.5
$ws.Range("I129").Value = 1113764.2
$ws.Range("J129").Value = 6539305.5
$ws.Range("K129").Value = 3341292.6
$ws.Range("L129").Value = 19617916.5
$ws.Range("M129").Value = -3336292.6
$ws.Range("N129").Value = -19627916.5
$ws.Range("H131").Value = 2925655
$ws.Range("J131").Value = 3474025.8
$ws.Range("L131").Value = 10422077.4
$ws.Range("N131").Value = -10432157.4
$ws.Range("H134").Value = 4832.5
$ws.Range("I134").Value = 3284.2856
$ws.Range("K134").Value = 9852.856800000001
$ws.Range("M134").Value = -4782.856800000001
$ws.Range("H139").Value = 3279.3103
$ws.Range("I139").Value = 2221.4285
$ws.Range("K139").Value = 6664.2855
$ws.Range("M139").Value = -1524.2855
$ws.Range("H141").Value = 9183.799999999999
$ws.Range("J141").Value = 15000
$ws.Range("L141").Value = 45000
$ws.Range("N141").Value = -55360

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 28489.21
$ws.Range("I70").Value = 22716.584
$ws.Range("J70").Value = 41477.625
$ws.Range("K70").Value = 22716.584
$ws.Range("L70").Value = 41477.625
$ws.Range("M70").Value = -22446.584
$ws.Range("N70").Value = -42017.625
$ws.Range("H73").Value = 28489.21
$ws.Range("I73").Value = 22716.584
$ws.Range("J73").Value = 41477.625
$ws.Range("K73").Value = 22716.584
$ws.Range("L73").Value = 41477.625
$ws.Range("M73").Value = -21780.584
$ws.Range("N73").Value = -43349.625
$ws.Range("H122").Value = 3188
$ws.Range("I122").Value = 3061.8518
$ws.Range("K122").Value = 9185.555399999999
$ws.Range("M122").Value = -6735.555399999999
$ws.Range("H132").Value = 6103.2
$ws.Range("I132").Value = 5550.825
$ws.Range("K132").Value = 16652.475
$ws.Range("M132").Value = -14122.475

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 12000
$ws.Range("J3").Value = 12000
$ws.Range("L3").Value = 12000
$ws.Range("N3").Value = -12224
$ws.Range("H7").Value = 3072.6667
$ws.Range("I7").Value = 3335.5715
$ws.Range("K7").Value = 3335.5715
$ws.Range("M7").Value = -3223.5715
$ws.Range("H15").Value = 12000
$ws.Range("J15").Value = 12000
$ws.Range("L15").Value = 12000
$ws.Range("N15").Value = -12340
$ws.Range("H22").Value = 2644.0688
$ws.Range("I22").Value = 1614.3846
$ws.Range("J22").Value = 3480.6875
$ws.Range("K22").Value = 1614.3846
$ws.Range("L22").Value = 3480.6875
$ws.Range("M22").Value = -1319.3846
$ws.Range("N22").Value = -4070.6875
$ws.Range("H27").Value = 2644.0688
$ws.Range("I27").Value = 1614.3846
$ws.Range("J27").Value = 3480.6875
$ws.Range("K27").Value = 1614.3846
$ws.Range("L27").Value = 3480.6875
$ws.Range("M27").Value = -1507.3846
$ws.Range("N27").Value = -3694.6875
$ws.Range("H61").Value = 5648.3687
$ws.Range("I61").Value = 5684.3887
$ws.Range("K61").Value = 5684.3887
$ws.Range("M61").Value = -5482.3887
$ws.Range("H100").Value = 3249.75
$ws.Range("I100").Value = 2666.3333
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 2666.3333
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -2125.3333
$ws.Range("N100").Value = -6082
$ws.Range("H113").Value = 5648.3687
$ws.Range("I113").Value = 5684.3887
$ws.Range("K113").Value = 5684.3887
$ws.Range("M113").Value = -3514.3887
$ws.Range("H126").Value = 3072.6667
$ws.Range("I126").Value = 3335.5715
$ws.Range("K126").Value = 10006.7145
$ws.Range("M126").Value = -7536.7145
$ws.Range("H132").Value = 3018.9707
$ws.Range("I132").Value = 2487.5386
$ws.Range("J132").Value = 4746.125
$ws.Range("K132").Value = 7462.6158
$ws.Range("L132").Value = 14238.375
$ws.Range("M132").Value = -4932.6158
$ws.Range("N132").Value = -19298.375
$ws.Range("H136").Value = 3169.652
$ws.Range("I136").Value = 1900.0952
$ws.Range("J136").Value = 16500
$ws.Range("K136").Value = 5700.2856
$ws.Range("L136").Value = 49500
$ws.Range("M136").Value = -3150.2856
$ws.Range("N136").Value = -54600

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 10010
$ws.Range("I20").Value = 10010
$ws.Range("K20").Value = 10010
$ws.Range("M20").Value = -9770
$ws.Range("H54").Value = 12700
$ws.Range("J54").Value = 14266.667
$ws.Range("L54").Value = 14266.667
$ws.Range("N54").Value = -15306.667
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H96").Value = 55558720
$ws.Range("I96").Value = 66669764
$ws.Range("J96").Value = 3488
$ws.Range("K96").Value = 66669764
$ws.Range("L96").Value = 3488
$ws.Range("M96").Value = -66668391
$ws.Range("N96").Value = -6234
$ws.Range("H122").Value = 92280.234
$ws.Range("I122").Value = 6982.8
$ws.Range("K122").Value = 20948.4
$ws.Range("M122").Value = -18498.4
$ws.Range("H132").Value = 1945.9354
$ws.Range("I132").Value = 1746.2693
$ws.Range("J132").Value = 2984.2
$ws.Range("K132").Value = 5238.8079
$ws.Range("L132").Value = 8952.599999999999
$ws.Range("M132").Value = -2708.8079
$ws.Range("N132").Value = -14012.6
$ws.Range("H136").Value = 787.80646
$ws.Range("I136").Value = 787.80646
$ws.Range("K136").Value = 2363.41938
$ws.Range("M136").Value = 186.5806199999997
